$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename task labels in column A (rows 8 down to 4, matching authoring order) ---
$ws.Range("A8").Value = "Task 7:Create Necessary DB Tables"
$ws.Range("A7").Value = "Task 6: DB:Schema"
$ws.Range("A6").Value = "Task 5: Establish DB Env"
$ws.Range("A5").Value = "Task 4: Login Page/Modal"
$ws.Range("A4").Value = "Task 3: Register Page/Modal"
# A9:A11 ("Task 8:", "Task 9:", "Task 10:") are unchanged.

# --- Style the existing "Complete" status cells (bold, accent6 theme green) ---
$ws.Range("B2:B3").Font.Bold = $true
$ws.Range("B2:B3").Font.ThemeColor = 10

# --- Row 5: "In Progress" status (bold, bright blue), notes, and last-updated date ---
$ws.Range("B5").Value = "In Progress"
$ws.Range("B5").Font.Bold = $true
$ws.Range("B5").Font.Color = 15773696

# --- Fill in the Status column with "Pending" for the remaining task rows plus extra blank rows ---
$ws.Range("B4").Value = "Pending"
$ws.Range("C5").Value = "polishing UI Design"

$ws.Range("B6").Value = "Pending"
$ws.Range("B7").Value = "Pending"
$ws.Range("B8").Value = "Pending"
$ws.Range("B9").Value = "Pending"
$ws.Range("B10").Value = "Pending"
$ws.Range("B11").Value = "Pending"
$ws.Range("B12").Value = "Pending"
$ws.Range("B13").Value = "Pending"
$ws.Range("B14").Value = "Pending"
$ws.Range("B15").Value = "Pending"
$ws.Range("B16").Value = "Pending"
$ws.Range("B17").Value = "Pending"

# --- Date Last Updated for the "In Progress" row ---
$ws.Range("D5").Value = 45324
$ws.Range("D5").NumberFormat = "mm-dd-yy"

# --- Update the active selection to D7 ---
$ws.Range("D7").Select()
